# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Column A narrowed from 45 to 44 characters.
# (ColumnWidth adds a fixed ~0.8333 padding vs. the stored OOXML <col> width,
# so back that constant out to land exactly on 44.)
$ws.Columns.Item(1).ColumnWidth = 43.166666666666664

# --- Bad Drivers table -------------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.90.2.1"
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 585
$ws.Range("D3").Value = 98.8

# Totals row mirrors the single Bad Drivers row above.
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = 585

# --- Good Drivers table --------------------------------------------------
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B12").Value = 56018
$ws.Range("D12").Value = 100
$ws.Range("E12").ClearContents()

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B13").Value = 34244
$ws.Range("D13").Value = 100
$ws.Range("E13").ClearContents()

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 442178
# D14 percentage unchanged (99.9); only the vintage date moves.
$ws.Range("E14").Value = "'2024-11-10"

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B15").Value = 77849
# D15 percentage unchanged (99.9); only the vintage date moves.
$ws.Range("E15").Value = "'2021-08-18"

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "'2020-08-05"

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "'2019-12-14"

# Rows 18-24 (older driver-vintage entries) have rolled off the report.
$ws.Range("A18:E24").ClearContents()
